$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffix = "_old"
$newSuffix = "_new"
$fv2410 = "_FV2410"
$fv2504 = "_FV2504"

for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value()
    if ($val -like "*$oldSuffix") {
        $base = $val.Substring(0, $val.Length - $oldSuffix.Length)
        $cell.Value = "$base$fv2410"
    } elseif ($val -like "*$newSuffix") {
        $base = $val.Substring(0, $val.Length - $newSuffix.Length)
        $cell.Value = "$base$fv2504"
    }
}

$listObj = $ws.ListObjects.Add(1, $ws.Range("A1:U65"), $null, 1)
$listObj.Name = "Table1"

$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
